# Removed backwards compatibility for multi-sheet bulk app translations
# The "Menus_and_forms" sheet used to have icon_filepath_en / audio_filepath_en /
# icon_filepath_fra / audio_filepath_fra columns. These are replaced with the
# image_en / audio_en / image_fra / audio_fra columns (and the now-unused
# strings are dropped from the workbook once nothing references them).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menus_and_forms")

$ws.Range("E1").Value = "image_en"
$ws.Range("F1").Value = "audio_en"
$ws.Range("G1").Value = "image_fra"
$ws.Range("H1").Value = "audio_fra"

$ws.Activate()
$ws.Range("H2").Select()
